$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '26.079.46'
Set-TextValue $ws.Range('E2') '  -1.01%  '
Set-TextValue $ws.Range('D3') '1.644.84'
Set-TextValue $ws.Range('E3') '  -1.30%  '
Set-TextValue $ws.Range('E4') '  -0.67%  '
Set-TextValue $ws.Range('D5') '217.73'
Set-TextValue $ws.Range('E5') '  -0.70%  '
Set-TextValue $ws.Range('D6') '0.5196'
Set-TextValue $ws.Range('E6') '  -2.74%  '
Set-TextValue $ws.Range('D8') '0.2616'
Set-TextValue $ws.Range('E8') '  -1.80%  '
Set-TextValue $ws.Range('D9') '0.06282'
Set-TextValue $ws.Range('E9') '  -1.77%  '
Set-TextValue $ws.Range('E10') '  -2.18%  '
Set-TextValue $ws.Range('D11') '0.07758'
Set-TextValue $ws.Range('E11') '  -1.19%  '
Set-TextValue $ws.Range('D12') '4.472'
Set-TextValue $ws.Range('E12') '  -2.07%  '
Set-TextValue $ws.Range('D13') '1.651.05'
Set-TextValue $ws.Range('E13') '  -0.93%  '
Set-TextValue $ws.Range('D14') '1.870.92'
Set-TextValue $ws.Range('E14') '  -1.28%  '
Set-TextValue $ws.Range('D15') '0.5572'
Set-TextValue $ws.Range('E15') '  +0.64%  '
Set-TextValue $ws.Range('D16') '0.0₅7995'
Set-TextValue $ws.Range('E16') '  -2.45%  '
Set-TextValue $ws.Range('D17') '64.78'
Set-TextValue $ws.Range('E17') '  -1.71%  '
Set-TextValue $ws.Range('D18') '26.077.65'
Set-TextValue $ws.Range('E18') '  -1.12%  '
Set-TextValue $ws.Range('D19') '1.004'
Set-TextValue $ws.Range('E19') '  -0.61%  '
Set-TextValue $ws.Range('D20') '4.639'
Set-TextValue $ws.Range('E20') '  -0.95%  '
Set-TextValue $ws.Range('D21') '192.65'
Set-TextValue $ws.Range('E21') '  -0.65%  '
Set-TextValue $ws.Range('D22') '10.10'
Set-TextValue $ws.Range('E22') '  -1.92%  '
Set-TextValue $ws.Range('D23') '5.957'
Set-TextValue $ws.Range('E23') '  -1.46%  '
Set-TextValue $ws.Range('D24') '1.006'
Set-TextValue $ws.Range('E24') '  -0.61%  '
Set-TextValue $ws.Range('D25') '146.59'
Set-TextValue $ws.Range('E25') '  +0.09%  '
Set-TextValue $ws.Range('D26') '0.1201'
Set-TextValue $ws.Range('E26') '  -2.40%  '
Set-TextValue $ws.Range('D27') '7.162'
Set-TextValue $ws.Range('E27') '  -0.73%  '
Set-TextValue $ws.Range('D28') '15.95'
Set-TextValue $ws.Range('E28') '  -1.06%  '
Set-TextValue $ws.Range('D29') '1.483'
Set-TextValue $ws.Range('E29') '  -1.14%  '
Set-TextValue $ws.Range('D30') '0.05618'
Set-TextValue $ws.Range('E30') '  -4.28%  '
Set-TextValue $ws.Range('D31') '1.263'
Set-TextValue $ws.Range('E31') '  -1.68%  '
Set-TextValue $ws.Range('D32') '3.460'
Set-TextValue $ws.Range('E32') '  -4.40%  '
Set-TextValue $ws.Range('D33') '3.358'
Set-TextValue $ws.Range('E33') '  +2.20%  '
Set-TextValue $ws.Range('D34') '1.598'
Set-TextValue $ws.Range('E34') '  -0.51%  '
Set-TextValue $ws.Range('D35') '2.795'
Set-TextValue $ws.Range('E35') '  -1.20%  '
Set-TextValue $ws.Range('D36') '2.411'
Set-TextValue $ws.Range('E36') '  -0.54%  '
Set-TextValue $ws.Range('D37') '0.9364'
Set-TextValue $ws.Range('E37') '  -3.46%  '
Set-TextValue $ws.Range('D38') '0.5664'
Set-TextValue $ws.Range('E38') '  -2.84%  '
Set-TextValue $ws.Range('D39') '5.958'
Set-TextValue $ws.Range('E39') '  +2.00%  '
Set-TextValue $ws.Range('D40') '0.01579'
Set-TextValue $ws.Range('E40') '  -1.42%  '
Set-TextValue $ws.Range('D41') '1.053.34'
Set-TextValue $ws.Range('E41') '  -1.07%  '
Set-TextValue $ws.Range('D42') '2.572'
Set-TextValue $ws.Range('E42') '  -0.50%  '
Set-TextValue $ws.Range('E43') '  -0.65%  '
Set-TextValue $ws.Range('D44') '0.8414'
Set-TextValue $ws.Range('E44') '  -2.53%  '
Set-TextValue $ws.Range('D45') '102.22'
Set-TextValue $ws.Range('E45') '  -2.32%  '
Set-TextValue $ws.Range('D46') '1.781.99'
Set-TextValue $ws.Range('E46') '  -1.34%  '
Set-TextValue $ws.Range('D47') '57.02'
Set-TextValue $ws.Range('E47') '  -1.49%  '
Set-TextValue $ws.Range('B48') 'Frax'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D48') '1.010'
Set-TextValue $ws.Range('E48') '  -0.30%  '
Set-TextValue $ws.Range('B49') 'BabyDogeCoin'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D49') '0.0₈105'
Set-TextValue $ws.Range('E49') '  -1.77%  '
Set-TextValue $ws.Range('D50') '0.05323'
Set-TextValue $ws.Range('E50') '  +3.00%  '
Set-TextValue $ws.Range('E51') '  -1.48%  '
